# Add team record (Wins/Losses/Ties) columns AD, AE, AF to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold/centered/bordered) from AC1 onto the three
# new header cells so AD1:AF1 inherit the same formatting as the other
# header row cells.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))

# Header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-48) gets the same team record values.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 69
    $ws.Cells.Item($r, 31).Value = 93
    $ws.Cells.Item($r, 32).Value = 0
}
